# Auto-generated Excel COM-interop script to apply crypto price/volume updates
# Commit: Updated cryptos list on Sat Jan 27 21:32:54 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Row, $Col, $Text)
    $cell = $ws.Cells.Item($Row, $Col)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = "Normal"
}

Set-TextCell 2 4 '42.098.81'  # D2
Set-TextCell 2 5 '  +0.40%  '  # E2
Set-TextCell 3 4 '2.272.56'  # D3
Set-TextCell 3 5 '  +0.81%  '  # E3
Set-TextCell 4 5 '  +0.06%  '  # E4
Set-TextCell 5 4 '305.85'  # D5
Set-TextCell 5 5 '  +1.30%  '  # E5
Set-TextCell 6 4 '93.44'  # D6
Set-TextCell 6 5 '  +1.44%  '  # E6
Set-TextCell 7 5 '  -0.25%  '  # E7
Set-TextCell 8 5 '  -0.01%  '  # E8
Set-TextCell 9 4 '0.490'  # D9
Set-TextCell 9 5 '  +1.48%  '  # E9
Set-TextCell 10 4 '32.87'  # D10
Set-TextCell 10 5 '  +1.78%  '  # E10
Set-TextCell 11 4 '0.0803'  # D11
Set-TextCell 11 5 '  +0.74%  '  # E11
Set-TextCell 12 5 '  -1.73%  '  # E12
Set-TextCell 13 4 '6.70'  # D13
Set-TextCell 13 5 '  +0.60%  '  # E13
Set-TextCell 14 4 '2.621.06'  # D14
Set-TextCell 14 5 '  +0.71%  '  # E14
Set-TextCell 15 4 '14.36'  # D15
Set-TextCell 15 5 '  +2.06%  '  # E15
Set-TextCell 16 4 '2.268.50'  # D16
Set-TextCell 16 5 '  +0.56%  '  # E16
Set-TextCell 17 5 '  +3.90%  '  # E17
Set-TextCell 18 4 '41.959.01'  # D18
Set-TextCell 18 5 '  +0.41%  '  # E18
Set-TextCell 19 4 '12.76'  # D19
Set-TextCell 19 5 '  +4.98%  '  # E19
Set-TextCell 20 4 '0.0₃0919'  # D20
Set-TextCell 20 5 '  +1.95%  '  # E20
Set-TextCell 21 5 '  +1.27%  '  # E21
Set-TextCell 22 4 '68.20'  # D22
Set-TextCell 22 5 '  +1.89%  '  # E22
Set-TextCell 23 4 '244.47'  # D23
Set-TextCell 23 5 '  +1.50%  '  # E23
Set-TextCell 24 2 'PancakeSwap'  # B24
Set-TextCell 24 3 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'  # C24
Set-TextCell 24 4 '2.61'  # D24
Set-TextCell 24 5 '  +2.31%  '  # E24
Set-TextCell 25 2 'ImmutableX'  # B25
Set-TextCell 25 3 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'  # C25
Set-TextCell 25 4 '1.95'  # D25
Set-TextCell 25 5 '  +2.99%  '  # E25
Set-TextCell 26 2 'Dai'  # B26
Set-TextCell 26 3 'https://coinranking.com/coin/MoTuySvg7+dai-dai'  # C26
Set-TextCell 26 4 '1.00'  # D26
Set-TextCell 26 5 '  +0.02%  '  # E26
Set-TextCell 27 2 'EthereumClassic'  # B27
Set-TextCell 27 3 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'  # C27
Set-TextCell 27 4 '24.06'  # D27
Set-TextCell 27 5 '  +0.38%  '  # E27
Set-TextCell 28 2 'Cosmos'  # B28
Set-TextCell 28 3 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'  # C28
Set-TextCell 28 4 '9.69'  # D28
Set-TextCell 28 5 '  +0.31%  '  # E28
Set-TextCell 29 2 'Toncoin'  # B29
Set-TextCell 29 3 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'  # C29
Set-TextCell 29 4 '2.09'  # D29
Set-TextCell 29 5 '  -9.19%  '  # E29
Set-TextCell 30 2 'BinanceUSD'  # B30
Set-TextCell 30 3 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'  # C30
Set-TextCell 30 4 '5.70'  # D30
Set-TextCell 30 5 '  +469.50%  '  # E30
Set-TextCell 31 4 '35.12'  # D31
Set-TextCell 31 5 '  +3.67%  '  # E31
Set-TextCell 32 4 '160.14'  # D32
Set-TextCell 32 5 '  +1.03%  '  # E32
Set-TextCell 33 5 '  +4.05%  '  # E33
Set-TextCell 35 4 '0.0745'  # D35
Set-TextCell 35 5 '  +0.31%  '  # E35
Set-TextCell 36 5 '  -0.43%  '  # E36
Set-TextCell 37 4 '17.16'  # D37
Set-TextCell 37 5 '  +4.41%  '  # E37
Set-TextCell 38 5 '  -1.05%  '  # E38
Set-TextCell 39 5 '  +1.51%  '  # E39
Set-TextCell 40 5 '  +1.08%  '  # E40
Set-TextCell 41 5 '  +0.72%  '  # E41
Set-TextCell 42 4 '4.01'  # D42
Set-TextCell 42 5 '  +2.21%  '  # E42
Set-TextCell 43 4 '19.79'  # D43
Set-TextCell 43 5 '  +1.30%  '  # E43
Set-TextCell 44 4 '2.015.21'  # D44
Set-TextCell 44 5 '  -1.71%  '  # E44
Set-TextCell 45 5 '  +9.55%  '  # E45
Set-TextCell 46 5 '  +1.52%  '  # E46
Set-TextCell 47 4 '10.25'  # D47
Set-TextCell 47 5 '  +1.90%  '  # E47
Set-TextCell 48 5 '  +1.79%  '  # E48
Set-TextCell 49 4 '53.37'  # D49
Set-TextCell 49 5 '  +3.57%  '  # E49
Set-TextCell 51 4 '72.63'  # D51
Set-TextCell 51 5 '  +2.93%  '  # E51
